# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D2").Value = '26.727.44'
$ws.Range("D2").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E2").Value = '  +1.31%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D3").Value = '1.725.32'
$ws.Range("D3").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E3").Value = '  +0.15%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D4").Value = '0.9979'
$ws.Range("D4").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D5").Value = '240.44'
$ws.Range("D5").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E5").Value = '  -0.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D6").Value = '0.9984'
$ws.Range("D6").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E6").Value = '  -0.12%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D7").Value = '0.4833'
$ws.Range("D7").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E7").Value = '  -0.99%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D8").Value = '0.2578'
$ws.Range("D8").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E8").Value = '  -0.27%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D9").Value = '0.06179'
$ws.Range("D9").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E9").Value = '  -0.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D10").Value = '1.720.09'
$ws.Range("D10").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E10").Value = '  -0.31%  '

# Row 11
$ws.Range("E11").Value = '  +2.74%  '

# Row 12
$ws.Range("E12").Value = '  -1.74%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D13").Value = '0.6037'
$ws.Range("D13").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E13").Value = '  +1.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D14").Value = '4.461'
$ws.Range("D14").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E14").Value = '  -1.54%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D15").Value = '76.95'
$ws.Range("D15").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E15").Value = '  -0.23%  '

# Row 16
$ws.Range("E16").Value = '  -0.15%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D17").Value = '26.552.09'
$ws.Range("D17").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E17").Value = '  +0.66%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D18").Value = '0.9981'
$ws.Range("D18").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E18").Value = '  -0.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D19").Value = '0.000007154'
$ws.Range("D19").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E19").Value = '  -0.59%  '

# Row 20
$ws.Range("E20").Value = '  +0.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D21").Value = '1.941.58'
$ws.Range("D21").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E21").Value = '  -0.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D22").Value = '4.404'
$ws.Range("D22").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E22").Value = '  -0.92%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D23").Value = '8.550'
$ws.Range("D23").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E23").Value = '  +0.68%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D24").Value = '5.045'
$ws.Range("D24").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E24").Value = '  -1.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D25").Value = '140.15'
$ws.Range("D25").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E25").Value = '  +1.56%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D26").Value = '15.24'
$ws.Range("D26").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("E27").Value = '  +3.09%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D28").Value = '106.56'
$ws.Range("D28").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E28").Value = '  -0.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D29").Value = '1.368'
$ws.Range("D29").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E29").Value = '  -2.21%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D30").Value = '4.013'
$ws.Range("D30").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E30").Value = '  +2.42%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D31").Value = '0.07912'
$ws.Range("D31").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E31").Value = '  -1.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D32").Value = '3.659'
$ws.Range("D32").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E32").Value = '  -0.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D33").Value = '0.04509'
$ws.Range("D33").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E33").Value = '  +0.22%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D34").Value = '2.594'
$ws.Range("D34").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E34").Value = '  -0.34%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D35").Value = '0.9978'
$ws.Range("D35").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E35").Value = '  +0.15%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D36").Value = '0.6182'
$ws.Range("D36").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E36").Value = '  -0.92%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D37").Value = '0.9317'
$ws.Range("D37").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E37").Value = '  +0.01%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D38").Value = '2.007'
$ws.Range("D38").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E38").Value = '  +2.35%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D39").Value = '2.443'
$ws.Range("D39").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E39").Value = '  +2.20%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D40").Value = '0.9978'
$ws.Range("D40").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E40").Value = '  -0.11%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D41").Value = '0.01494'
$ws.Range("D41").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E41").Value = '  +1.36%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D42").Value = '5.605'
$ws.Range("D42").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E42").Value = '  +2.51%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D43").Value = '99.74'
$ws.Range("D43").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E43").Value = '  -0.42%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D44").Value = '0.3822'
$ws.Range("D44").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E44").Value = '  -0.41%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D45").Value = '6.790'
$ws.Range("D45").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E45").Value = '  -0.92%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D46").Value = '0.1151'
$ws.Range("D46").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E46").Value = '  -0.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D47").Value = '0.05362'
$ws.Range("D47").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E47").Value = '  -0.06%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D48").Value = '7.936'
$ws.Range("D48").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E48").Value = '  +2.75%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D49").Value = '30.01'
$ws.Range("D49").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E49").Value = '  -0.13%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"   # keep numeric-looking price as literal text
$ws.Range("D50").Value = '1.241'
$ws.Range("D50").ClearFormats()       # drop the temporary text format, restore default style
$ws.Range("E50").Value = '  +1.07%  '

# Row 51
$ws.Range("E51").Value = '  +0.76%  '
